# Refresh the coin Price (D) / Volume(1h) (E) columns to the latest
# coinranking.com snapshot (GitHub Actions run, 2023-11-09 11:38 UTC).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, $text) {
    # Leading apostrophe forces Excel to keep numeric-looking text as text
    # (e.g. '14.70' would otherwise become the number 14.7).
    $range.Value = "'" + $text
    # Drop the auto-applied quote-prefix/text format so the cell keeps the
    # workbook's default (unstyled) look, matching the rest of the column.
    $range.ClearFormats()
}

$ws.Range('D2').Value = '36.715.14'
$ws.Range('E2').Value = '  +3.77%  '
$ws.Range('D3').Value = '1.906.61'
$ws.Range('E3').Value = '  +1.25%  '
Set-TextValue $ws.Range('D5') '248.76'
$ws.Range('E5').Value = '  +0.76%  '
Set-TextValue $ws.Range('D6') '0.694'
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('E7').Value = '  +0.01%  '
Set-TextValue $ws.Range('D8') '46.93'
$ws.Range('E8').Value = '  +8.26%  '
$ws.Range('E9').Value = '  +4.34%  '
Set-TextValue $ws.Range('D10') '57.59'
$ws.Range('E10').Value = '  +6.23%  '
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('E12').Value = '  +2.32%  '
Set-TextValue $ws.Range('D13') '14.70'
$ws.Range('E13').Value = '  +8.79%  '
$ws.Range('E14').Value = '  +4.79%  '
$ws.Range('D15').Value = '2.185.49'
$ws.Range('E15').Value = '  +1.33%  '
Set-TextValue $ws.Range('D16') '5.08'
$ws.Range('E16').Value = '  +2.47%  '
$ws.Range('D17').Value = '1.903.11'
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('D18').Value = '36.726.85'
$ws.Range('E18').Value = '  +3.87%  '
Set-TextValue $ws.Range('D19') '74.15'
$ws.Range('E19').Value = '  +0.95%  '
$ws.Range('D20').Value = '0.0₃0852'
$ws.Range('E20').Value = '  +2.90%  '
Set-TextValue $ws.Range('D21') '13.60'
$ws.Range('E21').Value = '  +5.96%  '
Set-TextValue $ws.Range('D22') '249.92'
$ws.Range('E22').Value = '  +2.14%  '
Set-TextValue $ws.Range('D23') '5.13'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('E24').Value = '  +0.04%  '
Set-TextValue $ws.Range('D25') '2.51'
$ws.Range('E25').Value = '  -3.44%  '
$ws.Range('E26').Value = '  +1.65%  '
Set-TextValue $ws.Range('D27') '166.65'
$ws.Range('E27').Value = '  +1.21%  '
Set-TextValue $ws.Range('D28') '8.75'
$ws.Range('E28').Value = '  +1.18%  '
$ws.Range('E29').Value = '  +1.68%  '
$ws.Range('E30').Value = '  -0.16%  '
Set-TextValue $ws.Range('D31') '4.61'
$ws.Range('E31').Value = '  +7.32%  '
Set-TextValue $ws.Range('D32') '0.0610'
$ws.Range('E32').Value = '  +1.94%  '
Set-TextValue $ws.Range('D33') '1.94'
$ws.Range('E33').Value = '  +3.36%  '
$ws.Range('E34').Value = '  +2.79%  '
$ws.Range('E35').Value = '  +0.05%  '
Set-TextValue $ws.Range('D36') '0.0870'
$ws.Range('E36').Value = '  +18.51%  '
Set-TextValue $ws.Range('D37') '18.69'
$ws.Range('E37').Value = '  +54.66%  '
$ws.Range('E38').Value = '  -0.67%  '
Set-TextValue $ws.Range('D39') '0.866'
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('E40').Value = '  +1.20%  '
Set-TextValue $ws.Range('D41') '104.11'
$ws.Range('E41').Value = '  +7.00%  '
$ws.Range('E42').Value = '  +4.05%  '
Set-TextValue $ws.Range('D43') '17.77'
$ws.Range('E43').Value = '  +2.65%  '
Set-TextValue $ws.Range('D44') '2.86'
$ws.Range('E44').Value = '  +19.67%  '
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('D46').Value = '1.346.12'
$ws.Range('E46').Value = '  +2.89%  '
Set-TextValue $ws.Range('D47') '2.37'
$ws.Range('E47').Value = '  -1.46%  '
Set-TextValue $ws.Range('D48') '0.0816'
$ws.Range('E48').Value = '  +0.41%  '
Set-TextValue $ws.Range('D49') '2.81'
$ws.Range('E49').Value = '  +2.68%  '
Set-TextValue $ws.Range('D50') '6.45'
$ws.Range('E50').Value = '  +1.94%  '
$ws.Range('D51').Value = '2.087.92'
$ws.Range('E51').Value = '  +1.33%  '
